# Apply the BOM column-C fix: copy column D's lookup values into column C
# (overwriting the placeholder "X" values), and introduce two new shared
# strings ("NEW" / "NEW2") for the two rows whose D value also changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("file_2")

# Row -> new value for columns C and D ("" means D keeps its existing value,
# and C should simply be set to that same value).
$updates = @(
    @{ Row = 4;  Value = "B" },
    @{ Row = 5;  Value = "C" },
    @{ Row = 8;  Value = "NEW" },
    @{ Row = 9;  Value = "B" },
    @{ Row = 10; Value = "C" },
    @{ Row = 11; Value = "F" },
    @{ Row = 14; Value = "A" },
    @{ Row = 15; Value = "NEW2" },
    @{ Row = 16; Value = "H" },
    @{ Row = 19; Value = "C" },
    @{ Row = 20; Value = "G" }
)

foreach ($u in $updates) {
    $row = $u.Row
    $value = $u.Value

    # Column C gets the value and switches from center- to left-aligned
    # (matching column D's existing style/font/border/wrap, just alignment).
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $value
    $cCell.HorizontalAlignment = -4131  # xlLeft

    # Column D only changes for rows 8 and 15 (new parts "NEW"/"NEW2");
    # the rest already matched.
    $ws.Cells.Item($row, 4).Value = $value
}

# Sheet view tweaks: clear the scrolled top-left cell and move the active
# selection to C9.
$ws.Activate()
$ws.Range("C9").Select()
